$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Part 1: Swap the order of the "FORMATION ACADEMIQUE" and
# "EXPERIENCE PROFESSIONNELLE" sections (commit: "reverse pro exp and
# education sections"). We move the whole "FORMATION ACADEMIQUE" block
# (its heading paragraph plus the four date-range paragraphs under it) so
# that it now sits right after the "EXPERIENCE PROFESSIONNELLE" section's
# content and right before "COMPETENCES TECHNIQUES", instead of sitting
# right after "PROFIL" and before "EXPERIENCE PROFESSIONNELLE".
# ---------------------------------------------------------------------------

function Find-ParagraphIndexByExactText($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text.Trim() -eq $text) {
            return $i
        }
    }
    return -1
}

$formationIdx = Find-ParagraphIndexByExactText $d "FORMATION ACADEMIQUE"
if ($formationIdx -eq -1) {
    throw "Could not find FORMATION ACADEMIQUE heading paragraph"
}

# The block to move is the heading plus the next 4 paragraphs (the four
# education date-range lines).
$blockStartIdx = $formationIdx
$blockEndIdx = $formationIdx + 4

$competencesIdx = Find-ParagraphIndexByExactText $d "COMPETENCES TECHNIQUES"
if ($competencesIdx -eq -1) {
    throw "Could not find COMPETENCES TECHNIQUES heading paragraph"
}

$target = $d.Paragraphs.Item($competencesIdx)
$insertPoint = $target.Range.Start

# Copy each paragraph of the block, one at a time (preserves each
# paragraph's own pPr - w:pBdr / w:spacing / w:jc - individually), to just
# before "COMPETENCES TECHNIQUES".
for ($idx = $blockStartIdx; $idx -le $blockEndIdx; $idx++) {
    $p = $d.Paragraphs.Item($idx)
    $len = $p.Range.End - $p.Range.Start
    $ft = $p.Range.FormattedText
    $destRange = $d.Range($insertPoint, $insertPoint)
    $destRange.FormattedText = $ft
    $insertPoint = $insertPoint + $len
}

# Now remove the original block (paragraph indices for the original block
# are unchanged, since we inserted the copies further along in the
# document).
$origStart = $d.Paragraphs.Item($blockStartIdx)
$origEnd = $d.Paragraphs.Item($blockEndIdx)
$origRange = $d.Range($origStart.Range.Start, $origEnd.Range.End)
$origRange.Delete()

# ---------------------------------------------------------------------------
# Part 2: Reorder the bullet lines inside "COMPETENCES TECHNIQUES".
# All five lines share identical paragraph formatting, so we simply swap
# which text occupies which paragraph slot (leaving paragraph marks/
# formatting untouched), going from:
#   MLOps / Visualisation / Bases de données / Langages / ML-AI
# to:
#   Bases de données / Visualisation / MLOps / ML-AI / Langages
# ---------------------------------------------------------------------------

$mlopsIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "MLOps :*") {
        $mlopsIdx = $i
        break
    }
}
if ($mlopsIdx -eq -1) {
    throw "Could not find MLOps paragraph"
}

$newOrder = @(
    "Bases de données : SQL, MongoDB, Neo4j, Redis",
    "Visualisation : tableau",
    "MLOps : Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit",
    "ML/AI : Scikit-Learn, Keras, Tensorflow, Pandas, pySpark, XGboost, OpenCV, Matplotlib, Seaborn",
    "Langages : python, matlab, c, c++"
)

for ($k = 0; $k -lt $newOrder.Length; $k++) {
    $p = $d.Paragraphs.Item($mlopsIdx + $k)
    $r = $p.Range
    # Exclude the trailing paragraph mark from the replaced range so the
    # paragraph's own formatting (pPr) is left untouched.
    $r.End = $r.End - 1
    $r.Text = $newOrder[$k]
}

Write-Output "Sections reordered and competences bullets reordered."
